$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 233.71428
$ws.Range("I6").Value = 237.2
$ws.Range("K6").Value = 711.5999999999999
$ws.Range("M6").Value = -599.5999999999999

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H62").Value = 2478.8
$ws.Range("I62").Value = 798.3333
$ws.Range("J62").Value = 4999.5
$ws.Range("K62").Value = 798.3333
$ws.Range("L62").Value = 4999.5
$ws.Range("M62").Value = -174.3333
$ws.Range("N62").Value = -6247.5

$ws.Range("H65").Value = 2478.8
$ws.Range("I65").Value = 798.3333
$ws.Range("J65").Value = 4999.5
$ws.Range("K65").Value = 3991.6665
$ws.Range("L65").Value = 24997.5
$ws.Range("M65").Value = -871.6665000000003
$ws.Range("N65").Value = -31237.5

$ws.Range("H129").Value = 1784719.2
$ws.Range("J129").Value = 4950
$ws.Range("L129").Value = 14850
$ws.Range("N129").Value = -24850

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2938.5
$ws.Range("I45").Value = 2446.5715
$ws.Range("K45").Value = 2446.5715
$ws.Range("M45").Value = -2069.5715

$ws.Range("H74").Value = 703.8461
$ws.Range("I74").Value = 639
$ws.Range("K74").Value = 639
$ws.Range("M74").Value = 235

$ws.Range("H77").Value = 703.8461
$ws.Range("I77").Value = 639
$ws.Range("K77").Value = 3195
$ws.Range("M77").Value = 1173

$ws.Range("H122").Value = 3224.3022
$ws.Range("I122").Value = 1907.9354
$ws.Range("J122").Value = 6624.9165
$ws.Range("K122").Value = 5723.8062
$ws.Range("L122").Value = 19874.7495
$ws.Range("M122").Value = -3273.8062
$ws.Range("N122").Value = -24774.7495

$ws.Range("H132").Value = 3912.9333
$ws.Range("I132").Value = 1740.4
$ws.Range("J132").Value = 4999.2
$ws.Range("K132").Value = 5221.200000000001
$ws.Range("L132").Value = 14997.6
$ws.Range("M132").Value = -2691.200000000001
$ws.Range("N132").Value = -20057.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7579623
$ws.Range("I94").Value = 3397.182
$ws.Range("K94").Value = 3397.182
$ws.Range("M94").Value = -2946.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 5152.0713
$ws.Range("I31").Value = 1239.4
$ws.Range("J31").Value = 7325.778
$ws.Range("K31").Value = 1239.4
$ws.Range("L31").Value = 7325.778
$ws.Range("M31").Value = -944.4000000000001
$ws.Range("N31").Value = -7915.778

$ws.Range("H33").Value = 1031
$ws.Range("I33").Value = 1031
$ws.Range("K33").Value = 1031
$ws.Range("M33").Value = -652

$ws.Range("H34").Value = 5152.0713
$ws.Range("I34").Value = 1239.4
$ws.Range("J34").Value = 7325.778
$ws.Range("K34").Value = 1239.4
$ws.Range("L34").Value = 7325.778
$ws.Range("M34").Value = -1037.4
$ws.Range("N34").Value = -7729.778

$ws.Range("H35").Value = 700
$ws.Range("I35").Value = 700
$ws.Range("K35").Value = 700
$ws.Range("M35").Value = -406

$ws.Range("H58").Value = 2557.6
$ws.Range("I58").Value = 1063
$ws.Range("J58").Value = 4799.5
$ws.Range("K58").Value = 1063
$ws.Range("L58").Value = 4799.5
$ws.Range("M58").Value = -860
$ws.Range("N58").Value = -5205.5

$ws.Range("H122").Value = 303676.28
$ws.Range("J122").Value = 5917.4375
$ws.Range("L122").Value = 17752.3125
$ws.Range("N122").Value = -22652.3125

$ws.Range("H136").Value = 2557.6
$ws.Range("I136").Value = 1063
$ws.Range("J136").Value = 4799.5
$ws.Range("K136").Value = 3189
$ws.Range("L136").Value = 14398.5
$ws.Range("M136").Value = -639
$ws.Range("N136").Value = -19498.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 42574708
$ws.Range("I4").Value = 42574708
$ws.Range("K4").Value = 127724124
$ws.Range("M4").Value = -127724012

$ws.Range("H18").Value = 998.1111
$ws.Range("I18").Value = 90
$ws.Range("J18").Value = 1724.6
$ws.Range("K18").Value = 270
$ws.Range("L18").Value = 5173.799999999999
$ws.Range("M18").Value = -101
$ws.Range("N18").Value = -5511.799999999999

$ws.Range("H117").Value = 3705.3333
$ws.Range("I117").Value = 674.5
$ws.Range("J117").Value = 4571.2856
$ws.Range("K117").Value = 2023.5
$ws.Range("L117").Value = 13713.8568
$ws.Range("M117").Value = 1418.5
$ws.Range("N117").Value = -20597.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 10526639
$ws.Range("J97").Value = 50000284
$ws.Range("L97").Value = 50000284
$ws.Range("N97").Value = -50001276

$ws.Range("H122").Value = 8341.1
$ws.Range("I122").Value = 2140.3333
$ws.Range("J122").Value = 10998.571
$ws.Range("K122").Value = 6420.999899999999
$ws.Range("L122").Value = 32995.713
$ws.Range("M122").Value = -3970.999899999999
$ws.Range("N122").Value = -37895.713

$ws.Range("H132").Value = 1580.0555
$ws.Range("I132").Value = 1596.5294
$ws.Range("J132").Value = 1300
$ws.Range("K132").Value = 4789.5882
$ws.Range("L132").Value = 3900
$ws.Range("M132").Value = -2259.5882
$ws.Range("N132").Value = -8960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 979.2
$ws.Range("I46").Value = 878.8
$ws.Range("J46").Value = 1079.6
$ws.Range("K46").Value = 878.8
$ws.Range("L46").Value = 1079.6
$ws.Range("M46").Value = -690.8
$ws.Range("N46").Value = -1455.6

$ws.Range("H61").Value = 3372.36
$ws.Range("I61").Value = 2818.1875
$ws.Range("K61").Value = 2818.1875
$ws.Range("M61").Value = -2616.1875

$ws.Range("H113").Value = 3372.36
$ws.Range("I113").Value = 2818.1875
$ws.Range("K113").Value = 2818.1875
$ws.Range("M113").Value = -648.1875

$ws.Range("H122").Value = 4966.5
$ws.Range("I122").Value = 1999
$ws.Range("J122").Value = 5560
$ws.Range("K122").Value = 5997
$ws.Range("L122").Value = 16680
$ws.Range("M122").Value = -3547
$ws.Range("N122").Value = -21580

$ws.Range("H132").Value = 4626.3477
$ws.Range("I132").Value = 2719
$ws.Range("J132").Value = 5643.6
$ws.Range("K132").Value = 8157
$ws.Range("L132").Value = 16930.8
$ws.Range("M132").Value = -5627
$ws.Range("N132").Value = -21990.8

$ws.Range("H136").Value = 5056.6523
$ws.Range("I136").Value = 2549
$ws.Range("J136").Value = 8316.6
$ws.Range("K136").Value = 7647
$ws.Range("L136").Value = 24949.8
$ws.Range("M136").Value = -5097
$ws.Range("N136").Value = -30049.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13891370
$ws.Range("I62").Value = 1767.3334
$ws.Range("J62").Value = 22225132
$ws.Range("K62").Value = 1767.3334
$ws.Range("L62").Value = 22225132
$ws.Range("M62").Value = -1143.3334
$ws.Range("N62").Value = -22226380

$ws.Range("H65").Value = 13891370
$ws.Range("I65").Value = 1767.3334
$ws.Range("J65").Value = 22225132
$ws.Range("K65").Value = 8836.666999999999
$ws.Range("L65").Value = 111125660
$ws.Range("M65").Value = -5716.666999999999
$ws.Range("N65").Value = -111131900

$ws.Range("H81").Value = 12825713
$ws.Range("I81").Value = 6190.143
$ws.Range("J81").Value = 66667708
$ws.Range("K81").Value = 12380.286
$ws.Range("L81").Value = 133335416
$ws.Range("M81").Value = -11319.286
$ws.Range("N81").Value = -133337538

$ws.Range("H84").Value = 12825713
$ws.Range("I84").Value = 6190.143
$ws.Range("J84").Value = 66667708
$ws.Range("K84").Value = 61901.43
$ws.Range("L84").Value = 666677080
$ws.Range("M84").Value = -56597.43
$ws.Range("N84").Value = -666687688

$ws.Range("H100").Value = 725.5
$ws.Range("I100").Value = 370.6
$ws.Range("K100").Value = 741.2
$ws.Range("M100").Value = -200.2

$ws.Range("H122").Value = 3416.2778
$ws.Range("I122").Value = 3292.8667
$ws.Range("J122").Value = 4033.3333
$ws.Range("K122").Value = 9878.6001
$ws.Range("L122").Value = 12099.9999
$ws.Range("M122").Value = -7428.6001
$ws.Range("N122").Value = -16999.9999
